$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Replace-Cell($row, $col, $oldText, $newText) {
    # Find.Execute ignores the scoping Range and searches/replaces across the
    # whole document, which is unsafe here because several of the new values
    # coincide with old values of other cells (and vice versa). Setting the
    # cell Range's Text directly is properly scoped to that cell only and
    # preserves the run formatting (font/size) already on the range.
    $rng = $t.Cell($row, $col).Range
    $rng.Text = $newText
}

Replace-Cell 1 1 "41×46=" "73×11="
Replace-Cell 1 2 "44×95=" "55×60="
Replace-Cell 1 3 "85×68=" "58×61="
Replace-Cell 1 4 "61×32=" "84×49="
Replace-Cell 1 5 "36×76=" "94×39="

Replace-Cell 5 1 "24×20=" "95×33="
Replace-Cell 5 2 "51×75=" "41×26="
Replace-Cell 5 3 "36×38=" "80×39="
Replace-Cell 5 4 "73×69=" "40×23="
Replace-Cell 5 5 "71×80=" "83×65="

Replace-Cell 10 1 "26×15=" "15×33="
Replace-Cell 10 2 "28×66=" "71×74="
Replace-Cell 10 3 "95×94=" "66×84="
Replace-Cell 10 4 "39×47=" "53×79="
Replace-Cell 10 5 "90×40=" "28×14="

Replace-Cell 15 1 "41×97=" "86×64="
Replace-Cell 15 2 "57×26=" "59×21="
Replace-Cell 15 3 "26×30=" "77×99="
Replace-Cell 15 4 "72×70=" "79×15="
Replace-Cell 15 5 "73×11=" "41×61="

Replace-Cell 20 1 "37×22=" "15×97="
Replace-Cell 20 2 "20×55=" "40×92="
Replace-Cell 20 3 "95×54=" "53×17="
Replace-Cell 20 4 "43×94=" "29×70="
Replace-Cell 20 5 "13×51=" "87×88="

Write-Host "All replacements complete"
